# The sheet used to be a "contact list" (First/Last Name, Company, Website,
# Mobile for Abhishek Yadav @ Hello Selenium). Turn it into a small
# "legal entity / KYC" form with 4 label/value columns instead of 6.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix up formatting first, while the donor cells still have the exact
# formats we want, by copying (not recomputing) their styles. This avoids
# creating brand-new font/fill/style-table entries for formats that already
# exist in the workbook.
#   A1 currently has the "center" header style; B1/C1/D1 already have the
#   "left" header style we want everywhere.
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)   # xlPasteFormats

#   A2 currently has the "quote-prefixed/centered" style left over from the
#   old "1" row-number column; B2/C2 already have the plain bordered style.
$ws.Range("B2").Copy()
$ws.Range("A2").PasteSpecial(-4122)

#   D2 will hold a numeric-looking value that must stay text; F2 already
#   carries the quote-prefixed text style used for such values.
$ws.Range("F2").Copy()
$ws.Range("D2").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# Remove the old hyperlink (Website column) before deleting its column.
$ws.Hyperlinks.Delete()

# --- New labels (row 1) ---------------------------------------------------
$ws.Range("A1").Value = "Full Name of Legal Entity"
$ws.Range("B1").Value = "Entity/Customer Type"
$ws.Range("C1").Value = "KYC (ANZ Only)"
$ws.Range("D1").Value = "Business Identifier Number"

# --- New values (row 2) ---------------------------------------------------
$ws.Range("A2").Value = "ABCD#12345"
$ws.Range("B2").Value = "Australian Company"
$ws.Range("C2").Value = "ANZ#123"
# Leading apostrophe keeps this numeric-looking value stored as text
# (quote-prefixed), just like the original sheet did for similar values.
$ws.Range("D2").Value = "'123456789"

# Drop the now-unused E and F columns (old Website/Mobile columns).
$ws.Range("E:F").Delete()

# --- Column widths to fit the new content --------------------------------
$ws.Columns("A").ColumnWidth = 22.5
$ws.Columns("B").ColumnWidth = 19.67
$ws.Columns("C").ColumnWidth = 18
$ws.Columns("D").ColumnWidth = 25.67

# --- View state ------------------------------------------------------------
$ws.Range("D10").Select()
